$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (Q4, Q5, Q6, Q7) which are no longer part of the table
$ws.Range("A6:G9").EntireRow.Delete()

# Update the recalculated forecast-error values for the remaining rows
$ws.Range("B2").Value = 0.394890212994997
$ws.Range("C2").Value = 0.4079149820070181
$ws.Range("D2").Value = 0.2433766280851557
$ws.Range("E2").Value = 0.4933321681029484
$ws.Range("F2").Value = 0.3068621425284427

$ws.Range("B3").Value = 0.1644430722176581
$ws.Range("C3").Value = 0.2400847867952391
$ws.Range("D3").Value = 0.09290660822557242
$ws.Range("E3").Value = 0.3048058533322029
$ws.Range("F3").Value = 0.2705243554966509
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = 0.17209220167078
$ws.Range("C4").Value = 0.1992204944551758
$ws.Range("D4").Value = 0.06134042005773921
$ws.Range("E4").Value = 0.2476699821491075
$ws.Range("F4").Value = 0.1951144100731962
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = 0.269664002978933
$ws.Range("C5").Value = 0.269664002978933
$ws.Range("D5").Value = 0.09111259386596202
$ws.Range("E5").Value = 0.3018486274044691
$ws.Range("F5").Value = 0.1918015608035556
$ws.Range("G5").Value = 2
